# Insert a new weekly price record for "Sandia" (Vega Modelo de Temuco)
# as row 291, shifting the existing rows 291-341 down to 292-342.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 291..341 down by one row.
$ws.Rows.Item(291).Insert()

# Fill the newly inserted row 291 with the new record.
$ws.Cells.Item(291, 1).Value = 10
$ws.Cells.Item(291, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(291, 3).Value = 'La Araucanía'
$ws.Cells.Item(291, 4).Value = 44522
$ws.Cells.Item(291, 5).Value = 9
$ws.Cells.Item(291, 6).Value = 100112028
$ws.Cells.Item(291, 7).Value = 'Sandia'
$ws.Cells.Item(291, 8).Value = 'Sin especificar'
$ws.Cells.Item(291, 9).Value = 'Primera'
$ws.Cells.Item(291, 10).Value = 300
$ws.Cells.Item(291, 11).Value = 950
$ws.Cells.Item(291, 12).Value = 950
$ws.Cells.Item(291, 13).Value = 950
$ws.Cells.Item(291, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(291, 15).Value = 'Perú'
$ws.Cells.Item(291, 16).Value = 950
$ws.Cells.Item(291, 17).Value = 1
$ws.Cells.Item(291, 18).Value = 'Hortaliza'
